$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.895.04"
$ws.Range("E2").Value = "'  +0.41%  "
$ws.Range("D3").Value = "'2.098.16"
$ws.Range("E3").Value = "'  +10.00%  "
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'252.44"
$ws.Range("E5").Value = "'  +1.52%  "
$ws.Range("D6").Value = "'0.660"
$ws.Range("E6").Value = "'  -5.21%  "
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("D8").Value = "'48.94"
$ws.Range("E8").Value = "'  +4.65%  "
$ws.Range("D9").Value = "'60.66"
$ws.Range("E9").Value = "'  +5.07%  "
$ws.Range("D10").Value = "'0.376"
$ws.Range("E10").Value = "'  +0.96%  "
$ws.Range("D11").Value = "'0.0744"
$ws.Range("E11").Value = "'  -1.44%  "
$ws.Range("D12").Value = "'0.109"
$ws.Range("E12").Value = "'  +9.15%  "
$ws.Range("D13").Value = "'14.91"
$ws.Range("E13").Value = "'  +1.69%  "
$ws.Range("D14").Value = "'2.402.69"
$ws.Range("D15").Value = "'0.837"
$ws.Range("E15").Value = "'  +3.40%  "
$ws.Range("D16").Value = "'2.126.91"
$ws.Range("E16").Value = "'  +11.52%  "
$ws.Range("D17").Value = "'5.15"
$ws.Range("E17").Value = "'  +1.54%  "
$ws.Range("D18").Value = "'36.715.33"
$ws.Range("E18").Value = "'  -0.10%  "
$ws.Range("D19").Value = "'73.18"
$ws.Range("E19").Value = "'  -1.32%  "
$ws.Range("D20").Value = "'0.0₃0821"
$ws.Range("E20").Value = "'  -3.49%  "
$ws.Range("D21").Value = "'13.28"
$ws.Range("E21").Value = "'  -2.26%  "
$ws.Range("D22").Value = "'240.93"
$ws.Range("E22").Value = "'  -3.67%  "
$ws.Range("E23").Value = "'  +4.24%  "
$ws.Range("E24").Value = "'  -0.11%  "
$ws.Range("D25").Value = "'2.53"
$ws.Range("E25").Value = "'  +1.21%  "
$ws.Range("D26").Value = "'170.76"
$ws.Range("E26").Value = "'  +2.45%  "
$ws.Range("D27").Value = "'9.42"
$ws.Range("E27").Value = "'  +7.93%  "
$ws.Range("D28").Value = "'20.95"
$ws.Range("E28").Value = "'  +12.45%  "
$ws.Range("D29").Value = "'1.99"
$ws.Range("E29").Value = "'  -9.26%  "
$ws.Range("B30").Value = "'Stellar"
$ws.Range("C30").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = "'  -4.42%  "
$ws.Range("B31").Value = "'Gas"
$ws.Range("C31").Value = "'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D31").Value = "'25.68"
$ws.Range("E31").Value = "'  +36.13%  "
$ws.Range("B32").Value = "'ImmutableX"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "'  +28.52%  "
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.51"
$ws.Range("E33").Value = "'  -2.56%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0614"
$ws.Range("E34").Value = "'  +0.34%  "
$ws.Range("D35").Value = "'0.0922"
$ws.Range("E35").Value = "'  +5.94%  "
$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = "'  +21.12%  "
$ws.Range("D38").Value = "'1.85"
$ws.Range("E38").Value = "'  -4.46%  "
$ws.Range("D39").Value = "'4.10"
$ws.Range("E39").Value = "'  -4.77%  "
$ws.Range("E40").Value = "'  -9.57%  "
$ws.Range("E41").Value = "'  -0.79%  "
$ws.Range("E42").Value = "'  +7.75%  "
$ws.Range("D43").Value = "'97.93"
$ws.Range("E43").Value = "'  -5.98%  "
$ws.Range("D44").Value = "'16.85"
$ws.Range("E44").Value = "'  -4.82%  "
$ws.Range("D45").Value = "'2.77"
$ws.Range("E45").Value = "'  -3.76%  "
$ws.Range("D46").Value = "'1.344.22"
$ws.Range("E46").Value = "'  -0.18%  "
$ws.Range("D47").Value = "'0.0852"
$ws.Range("E47").Value = "'  +4.26%  "
$ws.Range("D48").Value = "'7.11"
$ws.Range("E48").Value = "'  +10.43%  "
$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "'  +2.82%  "
$ws.Range("D50").Value = "'2.280.06"
$ws.Range("E50").Value = "'  +9.46%  "
$ws.Range("D51").Value = "'2.26"
$ws.Range("E51").Value = "'  -4.30%  "
